$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (row 2) field "YEAR OF BAPTISM" is renamed to "YEAR OF BIRTH".
$ws.Range("E2").Value = "YEAR OF BIRTH"

# Reflect that E2 is the last selected/edited cell.
$ws.Range("E2").Select()
